# Removed Test Case Inter-Dependency
# - rename the product's long name (productname) on both sheets
# - change shortname from a hardcoded numeric literal to a text placeholder
# - move the active selection off the input sheet's repayment-strategy row
#   and onto the output sheet (removes the prior test's leftover selection
#   state so this test case no longer depends on it)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoanInput
$ws2 = $wb.Worksheets.Item(2)   # ProductLoanOutput

$newProductName = "4299-MS-EI-DB-SAR-REC-RNI-FEE-FFConMONTHLYonLASTSUNDAY-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1st"

# productname (B1) updated on both the input and output sheets
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# shortname (B2) becomes the text "429c" instead of the number 4299
$ws1.Range("B2").Value = "429c"

# Move the selection on the input sheet from B17 to B21 and make the
# output sheet the active tab, so the workbook no longer re-opens onto
# the previous test's selection.
$ws1.Range("B21").Select()
$ws2.Activate()
